$d = $word.ActiveDocument

# 1. Fix LinkedIn URL to include www subdomain (for ATS auto-read)
$d.Content.Find.Execute(
    "https://linkedin.com/in/brian-phan-58530b1b0/",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://www.linkedin.com/in/brian-phan-58530b1b0/",
    2)

# 2. Merge "OSU COLLEGE OF AGRICULTURAL SCIENCES" split runs into one run
$d.Content.Find.Execute(
    "OSU COLLEGE OF AGRICULTURAL SCIENCES",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "OSU COLLEGE OF AGRICULTURAL SCIENCES",
    2)

# 3. Merge "9/202" + "5" into "9/2025"
$d.Content.Find.Execute(
    "9/202" + "5",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "9/2025",
    2)

# 4. Merge "Corvallis, United States - Software Developer" split runs into one
$d.Content.Find.Execute(
    "Corvallis, United States " + [char]8211 + " Software Developer",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Corvallis, United States " + [char]8211 + " Software Developer",
    2)

# 5. Merge "Maintained the portal..." split runs into one
$d.Content.Find.Execute(
    "Maintained the portal through which school districts and teachers can access the application, reporting systems, resources for developing outdoor school programs.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Maintained the portal through which school districts and teachers can access the application, reporting systems, resources for developing outdoor school programs.",
    2)

# 6. Merge "ity, Vietnam - Software Developer" split runs into one
$d.Content.Find.Execute(
    "ity, Vietnam " + [char]8211 + " Software Developer",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ity, Vietnam " + [char]8211 + " Software Developer",
    2)

# 7. Merge "Developed drone system for object detection, tracking and surveillance." split runs into one
$d.Content.Find.Execute(
    "Developed drone system for object detection, tracking and surveillance.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Developed drone system for object detection, tracking and surveillance.",
    2)
